$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width changes from 15.7109375 to 15.42578125 (matching column B).
# The runtime quantizes ColumnWidth to 1/6-character steps, so 14.6 is the
# input value whose stored width rounds to the closest achievable value (15.5).
$ws.Columns.Item(1).ColumnWidth = 14.6

$ws.Range("A1").Value = -0.04517650658601724
$ws.Range("B1").Value = 0.044741969204245891
$ws.Range("A2").Value = -0.022640681188736522
$ws.Range("B2").Value = 0.020824303131433908
$ws.Range("A3").Value = -0.014767646447232607
$ws.Range("B3").Value = 0.014475359423800072
$ws.Range("A4").Value = -0.0064753595048561152
$ws.Range("B4").Value = 0.006217954106388035
$ws.Range("A5").Value = -0.0032179541375114695
$ws.Range("B5").Value = 0.0023530804366140856
$ws.Range("A6").Value = -0.016367954639310511
$ws.Range("B6").Value = 0.016354117041732508
$ws.Range("A7").Value = -0.0063541171392071938
$ws.Range("B7").Value = 0.0063520800360894469
$ws.Range("A8").Value = 0.0036479198662058288
$ws.Range("B8").Value = -0.0036957919950784124
$ws.Range("A9").Value = 0.0056957919787192779
$ws.Range("B9").Value = -0.0057645928157392667
$ws.Range("A10").Value = -0.0091971876445029466
$ws.Range("B10").Value = 0.0091973688524564778
$ws.Range("A11").Value = -0.02437861539790287
$ws.Range("B11").Value = 0.024353703319777154
$ws.Range("A12").Value = -0.020853703349974495
$ws.Range("B12").Value = 0.020664720994435726
$ws.Range("A13").Value = -0.01716472102884925
$ws.Range("B13").Value = 0.017079033941992172
$ws.Range("A14").Value = -0.009079034021799437
$ws.Range("B14").Value = 0.0090514061368152099
$ws.Range("A15").Value = -0.008051406147587592
$ws.Range("B15").Value = 0.0080335996423288591
$ws.Range("A16").Value = -0.0060335996636164957
$ws.Range("B16").Value = 0.0060030838490496841
$ws.Range("A17").Value = -0.0040030838709324001
$ws.Range("B17").Value = 0.0039999999579771739
$ws.Range("A18").Value = -0.016101288075489606
$ws.Range("B18").Value = 0.016090841400568934
$ws.Range("A19").Value = -0.012090841440544065
$ws.Range("B19").Value = 0.0120159121392307
$ws.Range("A20").Value = -0.0080159121824365798
$ws.Range("B20").Value = 0.0080055533491201203
$ws.Range("A21").Value = -0.0040055533928544662
$ws.Range("B21").Value = 0.0039999999558339994
$ws.Range("A22").Value = -0.045715846617570222
$ws.Range("B22").Value = 0.045501785862100519
$ws.Range("A23").Value = -0.04050178591673248
$ws.Range("B23").Value = 0.040099417589376074
$ws.Range("A24").Value = -0.020099417799464447
$ws.Range("B24").Value = 0.019999999786656453
$ws.Range("A25").Value = 0.016014874183715833
$ws.Range("B25").Value = -0.016113582268921789
$ws.Range("A26").Value = -0.058574004994271789
$ws.Range("B26").Value = 0.058507101078005164
$ws.Range("A27").Value = -0.056007101103705992
$ws.Range("B27").Value = 0.055621404317725442
$ws.Range("A28").Value = -0.053621404348832336
$ws.Range("B28").Value = 0.053372240176157604
$ws.Range("A29").Value = -0.046372240263723441
$ws.Range("B29").Value = 0.0463117223654681
$ws.Range("A30").Value = 0.013688277017059658
$ws.Range("B30").Value = -0.013718710271332757
$ws.Range("A31").Value = 0.02071871018267224
$ws.Range("B31").Value = -0.02072725409015419
$ws.Range("A32").Value = -0.0040008293136999384
$ws.Range("B32").Value = 0.0039999999422253296
